# Update figures on the "PO List" sheet.
# Rows are keyed by the PO_name in column B:
#   Row 12 -> Brian Phua
#   Row 13 -> Lee Xuan Yen
#   Row 21 -> Santoso
#   Row 25 -> Jimmy Wong

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO List")

# Row 12 - Brian Phua
$ws.Range("H12").Value = 9
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 44900
$ws.Range("K12").Value = 2
$ws.Range("P12").Value = 4
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 44900
$ws.Range("S12").Value = 8
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 44900

# Row 13 - Lee Xuan Yen
$ws.Range("K13").Value = 3

# Row 21 - Santoso
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 44376
$ws.Range("K21").Value = 23

# Row 25 - Jimmy Wong
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 44902
$ws.Range("K25").Value = 1
